# "Generate Report for Handback" — mark the en-US -> zh-cn / de-de
# localization rows as handed back: update the Overview + per-language
# status text, stamp the handback file names / datetimes, and link the
# handed-back target file the same way the existing source-file link does.

$wb  = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns (E2, F2) ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# --- zh-cn sheet: Status column (C2) ---
$wsZhCn.Range("C2").Value = $newStatus

# --- de-de sheet: Status column (C2) ---
$wsDeDe.Range("C2").Value = $newStatus

# --- zh-cn sheet row 2: Latest Target File (I2) / Latest Handback File (J2) ---
$mdFileName  = "fed28c1a-55ce-41a9-9fd9-b5c3ed70f9f0.md"
$mdFileUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e51926f4884cc07c92e975be6a986ac1ddafdd7c/e2e/fed28c1a-55ce-41a9-9fd9-b5c3ed70f9f0.md"

$wsZhCn.Range("I2").Value = $mdFileName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdFileUrl, $null, $mdFileName, $mdFileName)
$wsZhCn.Range("I2").Font.Underline = $true
$wsZhCn.Range("I2").Font.Color = 15570276

$wsZhCn.Range("J2").Value = "fed28c1a-55ce-41a9-9fd9-b5c3ed70f9f0.271fe60cb78713f1a0a8bdfde164f2e7f81cd03d.zh-cn.xlf"

# Latest Handback DateTime (K2) - same slot as before, text refreshed
$wsZhCn.Range("K2").Value = "2016-08-16 14:59:50"

# --- de-de sheet row 2: Latest Target File (I2) / Latest Handback File (J2) ---
$wsDeDe.Range("I2").Value = $mdFileName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdFileUrl, $null, $mdFileName, $mdFileName)
$wsDeDe.Range("I2").Font.Underline = $true
$wsDeDe.Range("I2").Font.Color = 15570276

$wsDeDe.Range("J2").Value = "fed28c1a-55ce-41a9-9fd9-b5c3ed70f9f0.271fe60cb78713f1a0a8bdfde164f2e7f81cd03d.de-de.xlf"

# Latest Handback DateTime (K2) - new, later timestamp for de-de
$wsDeDe.Range("K2").Value = "2016-08-16 14:59:57"

# --- Column width refresh (status / target-file / handback-file columns
#     grew wider to fit the new text) ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.15
$wsOverview.Columns.Item(6).ColumnWidth = 29.15

$wsZhCn.Columns.Item(3).ColumnWidth  = 29.15
$wsZhCn.Columns.Item(9).ColumnWidth  = 39.15
$wsZhCn.Columns.Item(10).ColumnWidth = 39.15

$wsDeDe.Columns.Item(3).ColumnWidth  = 29.15
$wsDeDe.Columns.Item(9).ColumnWidth  = 39.15
$wsDeDe.Columns.Item(10).ColumnWidth = 39.15
